$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Update quantity for the 5V Single Channel Relay Module row (row 10): 1 -> 2
$ws.Range("B10").Value = 2

# Move the active selection to C13 on the Parts sheet
$ws.Range("C13").Select()
